# Refined metadata to be additional tab
#
# 1. Update the "time_taken" timestamps on the existing "data" sheet.
# 2. Add a new "metadata" worksheet (placed after "data") describing the
#    panel query that produced the data sheet.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh timestamps on the data sheet -------------------------------
$dataSheet.Range("F2").Value = "2021-10-05 14:33:42.131653"
$dataSheet.Range("F3").Value = "2021-10-05 14:33:42.131661"

# --- 2. Add the metadata sheet ----------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Reuse the header cell's existing format (bold, bordered, centered) for the
# header row and for A2, instead of building a brand-new style, so we don't
# introduce a duplicate entry in styles.xml.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Dent disease"
$metaSheet.Range("C2").Value = 96
# Keep "0.8" as text (not a number) to match the source data's data_version.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.8"
$metaSheet.Range("D2").Style = "Normal"
$metaSheet.Range("E2").Value = "2020-10-14T09:52:51.725872Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:42.128317"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/96/?format=json"

# Restore "data" as the active sheet/selection, matching the original
# workbook's active tab (unchanged by the diff).
$dataSheet.Activate()
[void]$dataSheet.Range("A1").Select()
